# Add data-source information to the header columns.
# The original sheet has two columns (B and C) that both contain the
# header text "FIRST_NAME" (column B coming from one data source, column
# C from a second data source). Disambiguate them by appending the
# source name to each header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "FIRST_NAME (Source1)"
$ws.Range("C1").Value = "FIRST_NAME (Source2)"
